$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = '62.245.67'
$ws.Range("E2").Value = '  -2.98%  '

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = '2.990.19'
$ws.Range("E3").Value = '  -3.98%  '

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.26%  '

# Row 5: BNB -> BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.93'
$ws.Range("E5").Value = '  -2.26%  '

# Row 6: Solana -> Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.86'
$ws.Range("E6").Value = '  -7.46%  '

# Row 7: USDC -> USDC
$ws.Range("E7").Value = '  +0.00%  '

# Row 8: XRP -> XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -4.10%  '

# Row 9: LidoStakedEther -> LidoStakedEther
$ws.Range("D9").Value = '2.988.76'
$ws.Range("E9").Value = '  -3.98%  '

# Row 10: Dogecoin -> Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -7.08%  '

# Row 11: Toncoin -> Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("E11").Value = '  -4.46%  '

# Row 12: Cardano -> Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -2.55%  '

# Row 13: ShibaInu -> ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("E13").Value = '  -5.42%  '

# Row 14: Avalanche -> Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.51'
$ws.Range("E14").Value = '  -7.16%  '

# Row 15: TRON -> TRON
$ws.Range("E15").Value = '  +1.63%  '

# Row 16: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '3.475.48'
$ws.Range("E16").Value = '  -4.30%  '

# Row 17: Polkadot -> Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.03'
$ws.Range("E17").Value = '  -2.98%  '

# Row 18: WrappedBTC -> WrappedBTC
$ws.Range("D18").Value = '62.217.89'
$ws.Range("E18").Value = '  -2.99%  '

# Row 19: WrappedEther -> WrappedEther
$ws.Range("D19").Value = '2.983.37'
$ws.Range("E19").Value = '  -4.31%  '

# Row 20: BitcoinCash -> BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '456.47'
$ws.Range("E20").Value = '  -4.67%  '

# Row 21: Chainlink -> Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.85'
$ws.Range("E21").Value = '  -4.49%  '

# Row 22: Polygon -> Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.677'
$ws.Range("E22").Value = '  -5.42%  '

# Row 23: Uniswap -> Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.28'
$ws.Range("E23").Value = '  -3.80%  '

# Row 24: Litecoin -> Fetch.AI
$ws.Range("B24").Value = 'Fetch.AI'
$ws.Range("C24").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("E24").Value = '  -7.90%  '

# Row 25: Fetch.AI -> Litecoin
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.83'
$ws.Range("E25").Value = '  -1.91%  '

# Row 26: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.23'
$ws.Range("E26").Value = '  -5.92%  '

# Row 27: Dai -> Dai
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.10%  '

# Row 28: RenderToken -> RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -5.97%  '

# Row 29: FirstDigitalUSD -> FirstDigitalUSD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.08%  '

# Row 30: NEARProtocol -> NEARProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  -4.36%  '

# Row 31: PancakeSwap -> PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.60'
$ws.Range("E31").Value = '  -3.62%  '

# Row 32: ImmutableX -> ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("E32").Value = '  -5.18%  '

# Row 33: EthereumClassic -> EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.81'
$ws.Range("E33").Value = '  -1.84%  '

# Row 34: Hedera -> Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -5.42%  '

# Row 35: Mantle -> Mantle
$ws.Range("E35").Value = '  -4.29%  '

# Row 36: PEPE -> PEPE
$ws.Range("D36").Value = '0.0₃0782'
$ws.Range("E36").Value = '  -7.45%  '

# Row 37: Filecoin -> Filecoin
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.72'
$ws.Range("E37").Value = '  -5.25%  '

# Row 38: Stacks -> Stacks
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.10'
$ws.Range("E38").Value = '  -7.11%  '

# Row 39: OKB -> OKB
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.94'
$ws.Range("E39").Value = '  -2.05%  '

# Row 40: Cosmos -> Cosmos
$ws.Range("E40").Value = '  -2.56%  '

# Row 41: dogwifhat -> dogwifhat
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("E41").Value = '  -11.92%  '

# Row 42: Bittensor -> Bittensor
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '408.94'
$ws.Range("E42").Value = '  -8.32%  '

# Row 43: TheGraph -> TheGraph
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.277'
$ws.Range("E43").Value = '  -5.51%  '

# Row 44: Kaspa -> Kaspa
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.111'
$ws.Range("E44").Value = '  -1.33%  '

# Row 45: Maker -> Maker
$ws.Range("D45").Value = '2.763.16'
$ws.Range("E45").Value = '  -2.45%  '

# Row 46: VeChain -> VeChain
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0349'
$ws.Range("E46").Value = '  -4.59%  '

# Row 47: Arweave -> Arweave
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.35'
$ws.Range("E47").Value = '  -6.61%  '

# Row 48: Monero -> Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.44'
$ws.Range("E48").Value = '  -2.56%  '

# Row 50: Stellar -> Stellar
$ws.Range("E50").Value = '  -2.75%  '

# Row 51: InjectiveProtocol -> InjectiveProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.68'
$ws.Range("E51").Value = '  -8.64%  '
